$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that are stored as literal text
# in the source workbook (e.g. "0.0001500" with significant trailing zeros).
# Force those cells to Text format before assignment so Excel does not
# auto-convert the string into a floating point number and lose precision
# / trailing zeros / distinguishing formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.24"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.75"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.401"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05685"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.327"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8053"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9185"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1404"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07451"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03125"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03030"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09376"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.874"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001583"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04802"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0005849"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006449"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004998"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001007"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.703"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.195"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1306"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04002"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006840"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1069"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002724"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007977"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005746"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4989"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2072"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
